# "updated everything following our discussion"
#
# - Remove the now-obsolete "Rhinorrhea" / "HP:0012373" / "Date of onset"
#   row-triplet (rows 5:7) from the "conditions" sheet. Deleting the whole
#   rows (rather than just clearing cells) shifts everything below them up,
#   which is what collapses the former rows 16:20 down to rows 13:17 and
#   shrinks the used range from A1:G20 to A1:G17.
# - Move the cursor/selection to A12 (matches the saved selection in the
#   workbook after the edit).
# - Set the page setup to Letter-ish defaults (paper size + portrait
#   orientation) as recorded in the saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("conditions")

$ws.Activate() | Out-Null

# Delete rows 5-7 entirely (shifts rows 8+ up by 3).
$ws.Rows("5:7").Delete() | Out-Null

# Move the active selection to A12.
$ws.Range("A12").Select() | Out-Null

# Record page setup (paper size 9 = A4, portrait orientation).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
